# Bold the term "DeepRanking" inside the bullet list of the "Subtitle 2"
# placeholder on slide 1 ("Adopted an implementation of DeepRanking to find
# similarity between images and built an image based query retrieval engine
# to get the top 5 best images."), splitting that run into three runs:
#   "Adopted an implementation of " | "DeepRanking" (bold) | " to find ..."

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item("Subtitle 2")
$tr = $sh.TextFrame.TextRange

$word = "DeepRanking"
$fullText = $tr.Text
$offset = $fullText.IndexOf($word)

if ($offset -ge 0) {
    $wordStart = $tr.Start + $offset
    $boldRange = $tr.Characters($wordStart, $word.Length)
    $boldRange.Font.Bold = 1
}
